# Supplementary Table 12 - intrinsic hydrolysis kinetics
# Reorganizing the table: update caption wording and switch the rate /
# std.dev. columns to a cleaner one-decimal scientific number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the caption text (A1) -----------------------------------
# The cell holds rich text: "Supplementary Table 12" (bold run) followed
# by a second run describing the figure. Only the wording about the
# number of replicates changes, so edit just that substring and keep the
# two runs' own formatting intact.
$cell = $ws.Range("A1")
$full = $cell.Characters().Text
$target = " three"
$idx = $full.IndexOf($target)
if ($idx -ge 0) {
    $startPos = $idx + 1
    $cell.Characters($startPos, $target.Length).Text = "2-9"
}

# Re-assert each run's own formatting so the rich-text split stays where
# it originally was (bold title vs. regular Arial description).
$newFull = $cell.Characters().Text
$titleLen = "Supplementary Table 12".Length
$descLen = $newFull.Length - $titleLen
$descRun = $cell.Characters($titleLen + 1, $descLen)
$descRun.Font.Name = "Arial"
$descRun.Font.Size = 12
$descRun.Font.Bold = $false

# --- 2. Give the rate / std. dev. columns a 1-decimal scientific format
$ws.Range("B3:C25").NumberFormat = "0.0E+00"

# --- 3. Update the saved selection/active cell --------------------------
$ws.Range("E9").Select() | Out-Null
